$wb = $excel.ActiveWorkbook

# 1. Rename the second worksheet
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Include #0"

# 2. Update Metadata sheet
$ws1 = $wb.Worksheets.Item(1)

# Version 0.1.0 -> 0.1.1
$ws1.Range("B3").Value = "0.1.1"

# Date update
$ws1.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# Insert a new row after "Contact" (row 10) for "Jurisdiction" with an empty value.
# Insert the row, then copy the formatting of the row above down onto it so the
# new row matches the rest of the table (border/fill/alignment), then set its values.
$ws1.Rows.Item(11).Insert()
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""
